# Fixed Bill bifurcation issue
# Populate column A (rows 2-7) with text values "1" through "6"
# using the column's existing text-formatted style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 1; $i -le 6; $i++) {
    $row = $i + 1
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = "$i"
    $cell.NumberFormat = "@"
}

# Move the active selection as seen in the edited workbook
$ws.Range("A11").Select()
